# Update cryptocurrency price (D) and 1h volume change (E) columns
# to the values from the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.322.28'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '2.761.52'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D5").Value = "'576.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = "'160.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("E10").Value = '  +4.87%  '
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '3.252.97'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = "'27.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").Value = '63.955.77'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = '2.769.12'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = "'12.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").Value = "'4.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").Value = "'358.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").Value = "'6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.09%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  -6.35%  '
$ws.Range("D24").Value = "'65.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").Value = "'8.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("E29").Value = '  +3.52%  '
$ws.Range("D30").Value = "'1.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("E31").Value = '  +9.93%  '
$ws.Range("D32").Value = "'168.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").Value = "'20.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.65%  '
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = "'1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").Value = "'350.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.34%  '
$ws.Range("D40").Value = "'6.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.85%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").Value = "'39.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").Value = "'22.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("D44").Value = "'21.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.10%  '
$ws.Range("D45").Value = "'0.0595"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("D46").Value = "'136.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").Value = "'0.632"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("D48").Value = "'0.0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").Value = '2.148.47'
$ws.Range("E51").Value = '  +1.45%  '
